$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("I3").Value = 3.5
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("W3").Value = 6.5
$ws.Range("AC3").Value = 8
$ws.Range("AI3").Value = 13
$ws.Range("AJ3").Value = 41
$ws.Range("AT3").Value = 2.5
$ws.Range("AV3").Value = 67
$ws.Range("AW3").Value = 5.5
$ws.Range("AX3").Value = 21

# Row 4
$ws.Range("G4").Value = 2.2
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 3.25
$ws.Range("J4").Value = 2.88
$ws.Range("L4").Value = 3.6
$ws.Range("Q4").Value = 1.85
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 1.36
$ws.Range("T4").Value = 3
$ws.Range("W4").Value = 9
$ws.Range("AB4").Value = 23
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 12
$ws.Range("AK4").Value = 23
$ws.Range("AL4").Value = 29
$ws.Range("AM4").Value = 151
$ws.Range("AS4").Value = 126
$ws.Range("AT4").Value = 3
$ws.Range("AW4").Value = 5
$ws.Range("AY4").Value = 23
$ws.Range("BA4").Value = 67

# Row 5
$ws.Range("G5").Value = 1.73
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 2.4
$ws.Range("L5").Value = 6
$ws.Range("Q5").Value = 2.35
$ws.Range("R5").Value = 1.57
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("X5").Value = 7
$ws.Range("Z5").Value = 13
$ws.Range("AA5").Value = 17
$ws.Range("AD5").Value = 7
$ws.Range("AG5").Value = 10
$ws.Range("AH5").Value = 23
$ws.Range("AN5").Value = 3.5
$ws.Range("AO5").Value = 9.5
$ws.Range("AQ5").Value = 34
$ws.Range("AV5").Value = 81
$ws.Range("AW5").Value = 6.5
$ws.Range("AZ5").Value = 126

# Row 6
$ws.Range("I6").Value = 2.4
$ws.Range("L6").Value = 3.4
$ws.Range("S6").Value = 1.62
$ws.Range("T6").Value = 2.2
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = 1.57
$ws.Range("AC6").Value = 6
$ws.Range("AK6").Value = 26
$ws.Range("AR6").Value = 126
$ws.Range("AT6").Value = 2.2
$ws.Range("AU6").Value = 9.5
$ws.Range("AY6").Value = 34
$ws.Range("BA6").Value = 101

# Row 7
$ws.Range("G7").Value = 2.25
$ws.Range("I7").Value = 3.1
$ws.Range("J7").Value = 3
$ws.Range("L7").Value = 3.75
$ws.Range("Q7").Value = 2.2
$ws.Range("R7").Value = 1.65
$ws.Range("U7").Value = 1.91
$ws.Range("V7").Value = 1.8
$ws.Range("W7").Value = 7
$ws.Range("X7").Value = 10
$ws.Range("AC7").Value = 8.5
$ws.Range("AE7").Value = 17
$ws.Range("AI7").Value = 12
$ws.Range("AJ7").Value = 34
$ws.Range("AL7").Value = 41
$ws.Range("AM7").Value = 351
$ws.Range("AU7").Value = 8.5
$ws.Range("AX7").Value = 19
$ws.Range("AY7").Value = 29

# Row 8
$ws.Range("G8").Value = 1.73
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 2.5
$ws.Range("K8").Value = 2.05
$ws.Range("L8").Value = 5.5
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 3
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.65
$ws.Range("W8").Value = 5.5
$ws.Range("Z8").Value = 13
$ws.Range("AD8").Value = 6.5
$ws.Range("AN8").Value = 3.6
$ws.Range("AQ8").Value = 34
$ws.Range("AU8").Value = 9.5
$ws.Range("AW8").Value = 6.5
$ws.Range("AX8").Value = 29
$ws.Range("BA8").Value = 151
